# Version 4 - Final Version
#
# Applies four text edits inside the "Problem Statement" and "Vision"
# paragraphs of the Elevator Pitch document:
#
#   1. "...hackers result in successful penetration..."
#        -> "...hackers allowing more and more successful penetration..."
#   2. "...successful penetrations of Computer Systems become more and
#       more common."
#        -> "...successful penetrations of Computer Systems."
#   3. "Should you accept my proposal, any potential weaknesses will be
#       identified"
#        -> "Should you accept my proposal, any potential vulnerabilities
#            will be identified"
#   4. "eliminate or significantly reduce those weaknesses."
#        -> "eliminate or significantly reduce those vulnerabilities."

$d = $word.ActiveDocument

# 1) "result in" -> "allowing more and more" (leading word before
#    "successful penetration...")
$found1 = $d.Content.Find.Execute(
    "result in successful penetration",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "allowing more and more successful penetration",
    2)
Write-Host "Step 1 (result in -> allowing more and more):" $found1

# 2) Drop the duplicated "become more and more common" tail, replacing it
#    with a full stop right after "Computer Systems".
$found2 = $d.Content.Find.Execute(
    "successful penetrations of Computer Systems become more and more common.",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "successful penetrations of Computer Systems.",
    2)
Write-Host "Step 2 (drop become more and more common tail):" $found2

# 3) "weaknesses" -> "vulnerabilities" in the "Should you accept..." sentence
$found3 = $d.Content.Find.Execute(
    "Should you accept my proposal, any potential weaknesses will be identified",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Should you accept my proposal, any potential vulnerabilities will be identified",
    2)
Write-Host "Step 3 (weaknesses -> vulnerabilities, identified sentence):" $found3

# 4) "weaknesses" -> "vulnerabilities" in the closing sentence of that paragraph
$found4 = $d.Content.Find.Execute(
    "eliminate or significantly reduce those weaknesses.",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "eliminate or significantly reduce those vulnerabilities.",
    2)
Write-Host "Step 4 (weaknesses -> vulnerabilities, eliminate sentence):" $found4
